# Penambahan informasi diskon di module ImportExportData golongan dan produk
#
# Insert a new "DISKON" column between the existing "SATUAN"/"HARGA..." block
# and "STOK ETALASE" on the "produk" sheet (column G), pushing the old
# G/H/I (STOK ETALASE / STOK GUDANG / MINIMAL STOK GUDANG) columns one to the
# right (H/I/J).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("produk")

# Insert a new column at G - this shifts the old G:I columns (and their
# column-width definitions) to H:J automatically.
$ws.Columns.Item(7).Insert()

# Header for the freshly inserted column. It inherits formatting (style,
# borders, fill) from the column that used to be there (now column H),
# matching the rest of the header row.
$ws.Cells.Item(1, 7).Value = "DISKON"

# Give the new column a sensible bestFit-like width.
$ws.Columns.Item(7).ColumnWidth = 7

# Match the recorded selection at the end of the edit.
$ws.Range("H6").Select()
